# Auto-generated edit script
$d = $word.ActiveDocument

# Update the date line (first paragraph)
$d.Content.Find.Execute("2025-03-18 Tuesday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-03-19 Wednesday", 2) | Out-Null

# Update each table cell by position (row-major, 20 rows x 5 cols)
$t = $d.Tables.Item(1)
$values = @(
    "9-0=",
    "91+7=",
    "55+16=",
    "68+16=",
    "65+17=",
    "72-58=",
    "52+35=",
    "74-23=",
    "55+6=",
    "44+32=",
    "95-37=",
    "4+26=",
    "45-17=",
    "87-81=",
    "49-37=",
    "25+5=",
    "60+25=",
    "22+8=",
    "72+12=",
    "28+40=",
    "32-1=",
    "41-4=",
    "12+87=",
    "53-12=",
    "65-59=",
    "84-13=",
    "2+72=",
    "97-57=",
    "72+4=",
    "80-24=",
    "19+25=",
    "0+49=",
    "2+11=",
    "44+1=",
    "73+15=",
    "85-69=",
    "45+45=",
    "23+41=",
    "15+24=",
    "3+22=",
    "59-8=",
    "49+40=",
    "70+20=",
    "82-41=",
    "75+8=",
    "15+13=",
    "84+3=",
    "37-32=",
    "83-3=",
    "72-65=",
    "77-52=",
    "0+25=",
    "8+61=",
    "84+0=",
    "58-53=",
    "23+61=",
    "4-4=",
    "46-0=",
    "76-73=",
    "31+40=",
    "39+9=",
    "66-3=",
    "67-26=",
    "49+4=",
    "47-46=",
    "89-4=",
    "37+20=",
    "11+1=",
    "77+5=",
    "22+62=",
    "44+44=",
    "30+41=",
    "8+59=",
    "3+13=",
    "47+14=",
    "45-5=",
    "26+16=",
    "7+2=",
    "68-52=",
    "81-32=",
    "84-47=",
    "21+39=",
    "62+0=",
    "29+67=",
    "54-2=",
    "30+13=",
    "67+4=",
    "42+10=",
    "39+23=",
    "31+5=",
    "78-51=",
    "43-10=",
    "14+84=",
    "38+49=",
    "26+46=",
    "83-61=",
    "68-21=",
    "60+30=",
    "84+5=",
    "57-7="
)

$rows = 20
$cols = 5
$idx = 0
for ($r = 1; $r -le $rows; $r++) {
    for ($c = 1; $c -le $cols; $c++) {
        $cell = $t.Cell($r, $c)
        $cell.Range.Text = $values[$idx]
        $idx++
    }
}

Write-Host "Done. Updated" $idx "cells."
